$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20
$ws.Range("A20").Value = 111964095
$ws.Range("Q20").Value = 572522
$ws.Range("R20").Value = 6635185

# Row 21
$ws.Range("A21").Value = 111964105
$ws.Range("Q21").Value = 572626
$ws.Range("R21").Value = 6635265

# Row 22
$ws.Range("A22").Value = 111964103
$ws.Range("Q22").Value = 572486
$ws.Range("R22").Value = 6635047

# Row 37
$ws.Range("A37").Value = 112103533
$ws.Range("Q37").Value = 572387
$ws.Range("R37").Value = 6635305

# Row 38
$ws.Range("A38").Value = 112103543
$ws.Range("Q38").Value = 572413
$ws.Range("R38").Value = 6635058

# Row 39
$ws.Range("A39").Value = 112103527
$ws.Range("B39").Value = 89096
$ws.Range("D39").Value = "NT"
$ws.Range("E39").Value = 5733
$ws.Range("F39").Value = "Såpfingersvamp"
$ws.Range("G39").Value = "Ramaria lutea"
$ws.Range("H39").Value = "(Vent.) Schild"
$ws.Range("K39").Value = $null
$ws.Range("Q39").Value = 572360
$ws.Range("R39").Value = 6635094

# Row 40
$ws.Range("A40").Value = 112103535
$ws.Range("B40").Value = 96735
$ws.Range("D40").Value = "VU"
$ws.Range("E40").Value = 220787
$ws.Range("F40").Value = "Knärot"
$ws.Range("G40").Value = "Goodyera repens"
$ws.Range("H40").Value = "(L.) R. Br."
$ws.Range("K40").Value = "fullt utvecklade blad"
$ws.Range("Q40").Value = 572400
$ws.Range("R40").Value = 6635245

# Row 57
$ws.Range("A57").Value = 111964085
$ws.Range("B57").Value = 89993
$ws.Range("D57").Value = "VU"
$ws.Range("E57").Value = 1209
$ws.Range("F57").Value = "Rynkskinn"
$ws.Range("G57").Value = "Phlebia centrifuga"
$ws.Range("H57").Value = "P.Karst."
# new (empty) placeholder cells that move in with the rotated row content
$ws.Range("J57").Font.Bold = $false
$ws.Range("K57").Font.Bold = $false
$ws.Range("N57").Font.Bold = $false
$ws.Range("P57").Value = "Rörbo, Sala-Norrby, Vstm"
$ws.Range("Q57").Value = 572348
$ws.Range("R57").Value = 6635252
$ws.Range("Y57").NumberFormat = "@"
$ws.Range("Y57").Value = "2023-09-08"
$ws.Range("AA57").NumberFormat = "@"
$ws.Range("AA57").Value = "2023-09-08"
$ws.Range("AF57").Font.Bold = $false

# Row 58
$ws.Range("A58").Value = 112103523
$ws.Range("B58").Value = 90826
$ws.Range("D58").Value = "LC"
$ws.Range("E58").Value = 4366
$ws.Range("F58").Value = "Skarp dropptaggsvamp"
$ws.Range("G58").Value = "Hydnellum peckii"
$ws.Range("H58").Value = "Banker"
$ws.Range("K58").Value = $null
$ws.Range("Q58").Value = 572413
$ws.Range("R58").Value = 6635058

# Row 59
$ws.Range("A59").Value = 112103539
$ws.Range("Q59").Value = 572445
$ws.Range("R59").Value = 6635165

# Row 60
$ws.Range("A60").Value = 112103547
$ws.Range("B60").Value = 96735
$ws.Range("E60").Value = 220787
$ws.Range("F60").Value = "Knärot"
$ws.Range("G60").Value = "Goodyera repens"
$ws.Range("H60").Value = "(L.) R. Br."
# J60 and N60 placeholder cells no longer present after rotation; K60 gains real text
$ws.Range("J60").Value = $null
$ws.Range("K60").Value = "fullt utvecklade blad"
$ws.Range("N60").Value = $null
$ws.Range("P60").Value = "Sala Norrby 1:3, Vstm"
$ws.Range("Q60").Value = 572369
$ws.Range("R60").Value = 6635135
$ws.Range("Y60").NumberFormat = "@"
$ws.Range("Y60").Value = "2023-09-13"
$ws.Range("AA60").NumberFormat = "@"
$ws.Range("AA60").Value = "2023-09-13"
$ws.Range("AF60").Value = $null

# Row 90
$ws.Range("A90").Value = 112276739
$ws.Range("B90").Value = 83506
$ws.Range("D90").Value = "NT"
$ws.Range("E90").Value = 241
$ws.Range("F90").Value = "Gransotdyna"
$ws.Range("G90").Value = "Camarops tubulina"
$ws.Range("H90").Value = "(Alb. & Schwein.:Fr.) Shear"
$ws.Range("K90").Value = $null
$ws.Range("Q90").Value = 572573
$ws.Range("R90").Value = 6635037

# Row 91
$ws.Range("A91").Value = 112276735
$ws.Range("B91").Value = 96735
$ws.Range("D91").Value = "VU"
$ws.Range("E91").Value = 220787
$ws.Range("F91").Value = "Knärot"
$ws.Range("G91").Value = "Goodyera repens"
$ws.Range("H91").Value = "(L.) R. Br."
$ws.Range("K91").Value = "fullt utvecklade blad"
$ws.Range("Q91").Value = 572075
$ws.Range("R91").Value = 6634976

# Row 92
$ws.Range("A92").Value = 112276718
$ws.Range("Q92").Value = 572500
$ws.Range("R92").Value = 6635075

# Row 93
$ws.Range("A93").Value = 112276716
$ws.Range("Q93").Value = 572567
$ws.Range("R93").Value = 6635215

# Row 95
$ws.Range("A95").Value = 112276743
$ws.Range("B95").Value = 90166
$ws.Range("D95").Value = "LC"
$ws.Range("E95").Value = 1339
$ws.Range("F95").Value = "Brandticka"
$ws.Range("G95").Value = "Pycnoporellus fulgens"
$ws.Range("H95").Value = "(Fr.) Donk"
$ws.Range("Q95").Value = 572657
$ws.Range("R95").Value = 6635285

# Row 96
$ws.Range("A96").Value = 112276692
$ws.Range("B96").Value = 94301
$ws.Range("D96").Value = "NT"
$ws.Range("E96").Value = 53
$ws.Range("F96").Value = "Vedtrappmossa"
$ws.Range("G96").Value = "Crossocalyx hellerianus"
$ws.Range("H96").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q96").Value = 572569
$ws.Range("R96").Value = 6635238
